$d = $word.ActiveDocument

# Locate the heading paragraph ("Suzanne and Joyce Wolfe Families"); the empty
# paragraph immediately before it currently hosts the "_GoBack" bookmark.
$headingPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Suzanne and Joyce Wolfe Families*") {
        $headingPara = $p
        break
    }
}

$headingStart = $headingPara.Range.Start

# Remove the paragraph mark of the preceding (empty) paragraph, merging it into
# the heading paragraph. This also drops the old "_GoBack" bookmark span (it was
# collapsed there), which we recreate below at its new location.
$mark = $d.Range($headingStart - 1, $headingStart)
$mark.Delete()

# After the merge, the heading paragraph starts right where the old empty
# paragraph used to start. Split "Suzanne and Joy" | "ce Wolfe Families" and
# drop the (now collapsed) bookmark between the two halves.
$newStart = $headingStart - 1
$splitPos = $newStart + "Suzanne and Joy".Length

$bmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
